$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("M2").Value = 8.979873666666666
$ws.Range("N2").Value = 26.939621
$ws.Range("O2").Value = 0.3651616045144693
$ws.Range("P2").Value = 0.3651616045144694
$ws.Range("Q2").Value = 0.3512956511312222
$ws.Range("R2").Value = 3.161660860181
$ws.Range("S2").Value = 0.3651616045144693
$ws.Range("T2").Value = 0.3651616045144694

# Row 3
$ws.Range("M3").Value = 4.482719
$ws.Range("N3").Value = 13.448157
$ws.Range("O3").Value = 0.1822872930499836
$ws.Range("P3").Value = 0.1822872930499837
$ws.Range("Q3").Value = 0.1753654615196667
$ws.Range("R3").Value = 1.578289153677
$ws.Range("S3").Value = 0.1822872930499836
$ws.Range("T3").Value = 0.1822872930499837

# Row 4
$ws.Range("M4").Value = 9.285498666666667
$ws.Range("N4").Value = 27.856496
$ws.Range("O4").Value = 0.3775896763919173
$ws.Range("P4").Value = 0.3775896763919173
$ws.Range("Q4").Value = 0.3632518030062222
$ws.Range("R4").Value = 3.269266227056
$ws.Range("S4").Value = 0.3775896763919173
$ws.Range("T4").Value = 0.3775896763919173

# Row 5
$ws.Range("M5").Value = 1.843414333333333
$ws.Range("N5").Value = 5.530243
$ws.Range("O5").Value = 0.07496142604362967
$ws.Range("P5").Value = 0.07496142604362969
$ws.Range("Q5").Value = 0.07211498319144445
$ws.Range("R5").Value = 0.649034848723
$ws.Range("S5").Value = 0.07496142604362967
$ws.Range("T5").Value = 0.07496142604362969
